$wb = $excel.ActiveWorkbook

# Rename existing sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "NewAccounts"

# Update header row to include AccountType/Type and Employees columns
$ws1.Range("C1").Value = "Type"
$ws1.Range("D1").Value = "Employees"

# Update data rows: AccountName values all become IntegrationTest, add Type & Employees
$accountName = "IntegrationTest"
$types = @("Prospect","Prospect","Channel Partner / Reseller","Installation Partner","Technology Partner","Other","Customer - Channel","Customer - Channel","Customer - Direct","Customer - Direct")
$employees = @(50,100,150,200,250,300,350,400,450,500)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws1.Range("A$row").Value = $accountName
    $ws1.Range("C$row").Value = $types[$i]
    $ws1.Range("D$row").Value = $employees[$i]
}

# Add second sheet after the first one
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "VisualforceToLWC"

$ws2.Range("A1").Value = "AccountName"
$ws2.Range("B1").Value = "AccountType"
$ws2.Range("C1").Value = "AccountPhone"
$ws2.Range("D1").Value = "Employees"

$ws2.Range("A2").Value = "Edge Communications"
$ws2.Range("B2").Value = "Customer - Direct"
$ws2.Range("C2").Value = 4158563255
$ws2.Range("D2").Value = 300

$ws2.Range("A3").Value = "Edge Communications"
$ws2.Range("B3").Value = "Customer - Channel"
$ws2.Range("C3").Value = 4155554323
$ws2.Range("D3").Value = 400

$ws2.Range("A4").Value = "Provar Webinar"
$ws2.Range("B4").Value = "Customer - Channel"
$ws2.Range("C4").Value = "555-555-5555"
$ws2.Range("D4").Value = 100

$ws2.Range("A5").Value = "Provar Webinar"
$ws2.Range("B5").Value = "Customer - Channel"
$ws2.Range("C5").Value = "555-555-5554"
$ws2.Range("D5").Value = 100

# Set selections to match the target workbook state
$ws2.Range("B2").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("C8").Select() | Out-Null
